$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interested count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 207
$ws1.Range("F4").Value = 823

# Sheet "全部类型" (All types) - same two events appear one row lower
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 207
$ws4.Range("F5").Value = 823
